$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1): E1:H1 -------------------------------------
# Typed in this order so new shared-string entries land at the same indices
# as the target workbook (G1 "lessHealth" was entered first, then E1, F1,
# and finally H1).
$ws.Range("G1").Value = "lessHealth"
$ws.Range("E1").Value = "minNumOfMoves"
$ws.Range("F1").Value = "giveChance"
$ws.Range("H1").Value = "with rounding"

# --- Column E: D/2 ----------------------------------------------------------
# E2 is entered on its own, then E3:E13 filled down from E3 (matches the
# shared-formula grouping recorded in the target workbook: E2 individual,
# E3:E13 shared).
$ws.Range("E2").Formula = "=D2/2"
$ws.Range("E3:E13").Formula = "=D3/2"

# --- Column G: 100/D ---------------------------------------------------------
$ws.Range("G2").Formula = "=(100/D2)"
$ws.Range("G3:G13").Formula = "=(100/D3)"

# --- Column H: ROUNDUP(G,0) --------------------------------------------------
# H2 and H3 are entered individually (not shared yet), then H4:H13 is filled
# down from H4, forming its own shared-formula group.
$ws.Range("H2").Formula = "=ROUNDUP(G2,0)"
$ws.Range("H3").Formula = "=ROUNDUP(G3,0)"
$ws.Range("H4:H13").Formula = "=ROUNDUP(G4,0)"

# --- Column widths -----------------------------------------------------------
$ws.Columns("C:D").ColumnWidth = 11.166666666666666
$ws.Columns("E").ColumnWidth = 16.0
$ws.Columns("F").ColumnWidth = 10.333333333333334

# --- Selection ----------------------------------------------------------------
$ws.Range("J11").Select()
